# Seed questions and answers with xlsx
# Adds a new "Answers" worksheet (after "Questions"), fills it with the
# answers data, turns the range into a table ("Tableau2"), and leaves the
# new sheet as the active/selected tab (matching the authored workbook).

$wb = $excel.ActiveWorkbook
$questions = $wb.Worksheets.Item(1)

# Insert the new sheet right after "Questions".
$answers = $wb.Worksheets.Add($null, $questions)
$answers.Name = "Answers"

# Header row.
$answers.Range("A1").Value = "text"
$answers.Range("B1").Value = "correct"
$answers.Range("C1").Value = "question_id"

# Data rows.
$answers.Range("A2").Value = "hdkjfhaf"
$answers.Range("B2").Value = 0
$answers.Range("C2").Value = 1

$answers.Range("A3").Value = "ajfdljasfdlkjdsf"
$answers.Range("B3").Value = 1
$answers.Range("C3").Value = 1

$answers.Range("A4").Value = "dsjflkjdsafljsalkfd"
$answers.Range("B4").Value = 0
$answers.Range("C4").Value = 1

$answers.Range("A5").Value = "jaslfdj"
$answers.Range("B5").Value = 0
$answers.Range("C5").Value = 1

# Column C is a touch wider than default, same as the authored sheet.
$answers.Columns.Item(3).ColumnWidth = 11.65

# Turn the filled range into a table, like "Tableau1" on the Questions sheet.
$table = $answers.ListObjects.Add(1, $answers.Range("A1:C5"), $null, 1)
$table.Name = "Tableau2"

# Leave the cursor parked below the table and make "Answers" the active tab.
$answers.Range("C6").Select() | Out-Null
